$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$ws.Range("A25").Value = "Deep Learning"

$ws.Range("C24").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = (Get-Date -Year 2020 -Month 2 -Day 13).Date

$ws.Range("E25").Value = "deep learning;machine learning;data science;neural networks"
$ws.Range("F25").Value = "Ebook"

[void]$ws.Range("G25").Select()
